$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the AutoFilter over the full data range (including the newly
# added row 36) and restrict column D (GeogAreaCode) to the "Desert Days"
# family of values, mirroring a refresh/sync from the source data feed.
$ws.AutoFilterMode = $false

$filterRange = $ws.Range("A1:E36")
$criteria = @("3N DESERT DAYS", "4N DESERT DAYS", "DESERT DAYS", "DESERTDAYS26", "DESERTDAYSDUBAI", "DESERTDAYSDUBAI26")
[void]$filterRange.AutoFilter(4, $criteria, 7)

# Keep the selection where the user last left it after refreshing.
[void]$ws.Range("D22").Select()
